# Week 11 Newsletter checkin: refresh injury report rows with latest
# updates, drop players who are off the DL (Tom Murphy, Adam Ottavino),
# and add the two new entries (Gerardo Parra, Chad Qualls).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Tyler Anderson - updated date + details
$ws.Range("C2").Value = "June 12 2017"
$ws.Range("E2").Value = "Anderson has been placed on the 10-day disabled list with knee inflammation but is likely to return to the starting rotation before the end of June."

# Row 3: Chad Bettis - updated date + details
$ws.Range("C3").Value = "June 06 2017"
$ws.Range("E3").Value = "Bettis is on the 60-day disabled list while recovering from testicular cancer but is expected to make his season debut sometime around the All-Star break."

# Row 4: David Dahl - updated date only
$ws.Range("C4").Value = "June 11 2017"

# Row 5: Jon Gray - updated date + details
$ws.Range("C5").Value = "June 16 2017"
$ws.Range("E5").Value = "Gray was placed on the 10-day disabled list with a stress fracture in his left foot. He is on a rehab assignment and is expected to rejoin the rotation during the end of June."

# Row 6: replace Tom Murphy with Gerardo Parra
$ws.Range("A6").Value = "Gerardo Parra"
$ws.Range("B6").Value = "parrage01"
$ws.Range("C6").Value = "June 07 2017"
$ws.Range("D6").Value = "Quadricep"
$ws.Range("E6").Value = "Parra has landed on the 10-day disabled list with a strained right quadriceps and is likely to be sidelined until the end of June."

# Row 7: replace Adam Ottavino with Chad Qualls
$ws.Range("A7").Value = "Chad Qualls"
$ws.Range("B7").Value = "quallch01"
$ws.Range("C7").Value = "June 18 2017"
$ws.Range("D7").Value = "Back"
$ws.Range("E7").Value = "Qualls has been placed on the 10-day disabled list with lower back spasms and it is unclear how much time he is `nexpected to miss."

$ws.Range("E8").Select() | Out-Null
